$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Uvar")
$ws2 = $wb.Worksheets.Item("UShock")

# --- Restructure the "Uvar" table -----------------------------------------
# Old layout (B..J): Meta, DYs, Dpae, t, g, DTI, Dps, Dpms, iext
# New layout (B..I):       DYs, ys, Dpae, t, g, DTI, Dps, iext
# 1) Drop the "Meta" column (old column B).
$ws1.Columns.Item(2).Delete()
# 2) Drop the "Dpms" column. After step 1 it shifted from column I to column H.
$ws1.Columns.Item(8).Delete()
# 3) Make room for the new "ys" column right after "DYs" (new column C).
$ws1.Columns.Item(3).Insert()

# Headers for the freshly inserted column.
$ws1.Cells.Item(1, 3).Value() = "ys"
$ws1.Cells.Item(2, 3).Value() = "res_ys"

# New data values for the "ys" column (rows 3-10).
# (written in plain decimal form -- this interpreter's numeric literal
# grammar does not accept scientific notation such as "1E-2")
$ysValues = @(
    -0.3051336351780229,
    -0.1715485635722717,
    -0.08826480561084699,
    -0.04179821773652505,
    -0.02040548508232441,
    -0.01501057065585601,
    -0.01901486379145872,
    -0.02785400560906104
)

for ($i = 0; $i -lt $ysValues.Count; $i++) {
    $row = 3 + $i
    $ws1.Cells.Item($row, 3).Value() = $ysValues[$i]
}

# --- Selection bookkeeping ---------------------------------------------------
[void]$ws2.Range("I13").Select()
[void]$ws1.Range("F16").Select()
